$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.787.78"
$ws.Range("E2").Value = "  +3.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.638.89"
$ws.Range("E3").Value = "  +7.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.27"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.28"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.623.12"
$ws.Range("E7").Value = "  +6.82%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.202"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.77"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "686.61"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.219.68"
$ws.Range("E15").Value = "  +6.82%  "
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.927.54"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.580.37"
$ws.Range("E18").Value = "  +5.03%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.34"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.935"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +8.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.78"
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.51"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +3.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.05"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.19"
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.20"
$ws.Range("E32").Value = "  +15.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "580.82"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.34"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.40"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.678.73"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.76"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0764"
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0469"
$ws.Range("E43").Value = "  +9.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.348"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.39"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.81"
$ws.Range("E47").Value = "  +5.14%  "
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.43"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.02"
$ws.Range("E51").Value = "  +0.05%  "
